# "Updated Documents with most recent version"
#
# Adds the newest team member, Steven Sewell (Community Manager, based in
# Las Vegas, NV), as a new row on the "Current" roster sheet, matching the
# formatting of the existing rows and widening the Email column so the
# longer address still reads cleanly.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Current")

# Append the new teammate as row 8 (right after Rob Solomon / row 7)
$ws1.Range("A8").Value = "Steven Sewell"
$ws1.Range("B8").Value = "stevenhasspam@gmail.com"
$ws1.Range("C8").Value = "714.552.9943"
$ws1.Range("D8").Value = "Community Manager"
$ws1.Range("E8").Value = "Las Vegas, NV"

# Re-use the same look as the rest of the data rows
$ws1.Range("A8:E8").Style = $ws1.Range("A7:E7").Style

# Widen the Email column a bit to comfortably fit the new address
$ws1.Columns("B").ColumnWidth = 26.0
